$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix R1487 and R1488: previously empty inline-string cells, now numeric 0
$ws.Cells.Item(1487, 18).Value = 0
$ws.Cells.Item(1488, 18).Value = 0

# Reuse the date/time number format already used by column A (style index 2)
$dateFmt = $ws.Cells.Item(1488, 1).NumberFormat

# Append 20 new weekly rows (1489-1508) carried over from the source export

# Row 1489
$ws.Cells.Item(1489, 1).NumberFormat = $dateFmt
$ws.Cells.Item(1489, 1).Value = 45474
$ws.Cells.Item(1489, 2).Value = 1484.949951171875
$ws.Cells.Item(1489, 3).Value = 1516
$ws.Cells.Item(1489, 4).Value = 1467
$ws.Cells.Item(1489, 5).Value = 1509.900024414062
$ws.Cells.Item(1489, 6).Value = 1497.205200195312
$ws.Cells.Item(1489, 7).Value = 7253948
$ws.Cells.Item(1489, 8).Value = 2024
$ws.Cells.Item(1489, 9).Value = 7
$ws.Cells.Item(1489, 10).Value = 1
$ws.Cells.Item(1489, 11).Value = 0
$ws.Cells.Item(1489, 12).Value = 0
$ws.Cells.Item(1489, 13).Value = 0
$ws.Cells.Item(1489, 14).Value = 27
$ws.Cells.Item(1489, 15).Value = 0
$ws.Cells.Item(1489, 16).Value = 0
$ws.Cells.Item(1489, 17).Value = 0

# Row 1490
$ws.Cells.Item(1490, 1).NumberFormat = $dateFmt
$ws.Cells.Item(1490, 1).Value = 45481
$ws.Cells.Item(1490, 2).Value = 1516
$ws.Cells.Item(1490, 3).Value = 1524.650024414062
$ws.Cells.Item(1490, 4).Value = 1483.199951171875
$ws.Cells.Item(1490, 5).Value = 1512.050048828125
$ws.Cells.Item(1490, 6).Value = 1499.337158203125
$ws.Cells.Item(1490, 7).Value = 6493934
$ws.Cells.Item(1490, 8).Value = 2024
$ws.Cells.Item(1490, 9).Value = 7
$ws.Cells.Item(1490, 10).Value = 8
$ws.Cells.Item(1490, 11).Value = 0
$ws.Cells.Item(1490, 12).Value = 0
$ws.Cells.Item(1490, 13).Value = 0
$ws.Cells.Item(1490, 14).Value = 28
$ws.Cells.Item(1490, 15).Value = 0
$ws.Cells.Item(1490, 16).Value = 0
$ws.Cells.Item(1490, 17).Value = 0

# Row 1491
$ws.Cells.Item(1491, 1).NumberFormat = $dateFmt
$ws.Cells.Item(1491, 1).Value = 45488
$ws.Cells.Item(1491, 2).Value = 1516.599975585938
$ws.Cells.Item(1491, 3).Value = 1531.949951171875
$ws.Cells.Item(1491, 4).Value = 1480.199951171875
$ws.Cells.Item(1491, 5).Value = 1485.5
$ws.Cells.Item(1491, 6).Value = 1473.010375976562
$ws.Cells.Item(1491, 7).Value = 5826714
$ws.Cells.Item(1491, 8).Value = 2024
$ws.Cells.Item(1491, 9).Value = 7
$ws.Cells.Item(1491, 10).Value = 15
$ws.Cells.Item(1491, 11).Value = 0
$ws.Cells.Item(1491, 12).Value = 0
$ws.Cells.Item(1491, 13).Value = 0
$ws.Cells.Item(1491, 14).Value = 29
$ws.Cells.Item(1491, 15).Value = 0
$ws.Cells.Item(1491, 16).Value = 0
$ws.Cells.Item(1491, 17).Value = 0

# Row 1492
$ws.Cells.Item(1492, 1).NumberFormat = $dateFmt
$ws.Cells.Item(1492, 1).Value = 45495
$ws.Cells.Item(1492, 2).Value = 1475.050048828125
$ws.Cells.Item(1492, 3).Value = 1600
$ws.Cells.Item(1492, 4).Value = 1470.650024414062
$ws.Cells.Item(1492, 5).Value = 1575
$ws.Cells.Item(1492, 6).Value = 1561.7578125
$ws.Cells.Item(1492, 7).Value = 11725125
$ws.Cells.Item(1492, 8).Value = 2024
$ws.Cells.Item(1492, 9).Value = 7
$ws.Cells.Item(1492, 10).Value = 22
$ws.Cells.Item(1492, 11).Value = 0
$ws.Cells.Item(1492, 12).Value = 0
$ws.Cells.Item(1492, 13).Value = 0
$ws.Cells.Item(1492, 14).Value = 30
$ws.Cells.Item(1492, 15).Value = 0
$ws.Cells.Item(1492, 16).Value = 0
$ws.Cells.Item(1492, 17).Value = 0

# Row 1493
$ws.Cells.Item(1493, 1).NumberFormat = $dateFmt
$ws.Cells.Item(1493, 1).Value = 45502
$ws.Cells.Item(1493, 2).Value = 1588
$ws.Cells.Item(1493, 3).Value = 1589
$ws.Cells.Item(1493, 4).Value = 1523.550048828125
$ws.Cells.Item(1493, 5).Value = 1528.800048828125
$ws.Cells.Item(1493, 6).Value = 1515.9462890625
$ws.Cells.Item(1493, 7).Value = 11500130
$ws.Cells.Item(1493, 8).Value = 2024
$ws.Cells.Item(1493, 9).Value = 7
$ws.Cells.Item(1493, 10).Value = 29
$ws.Cells.Item(1493, 11).Value = 0
$ws.Cells.Item(1493, 12).Value = 0
$ws.Cells.Item(1493, 13).Value = 0
$ws.Cells.Item(1493, 14).Value = 31
$ws.Cells.Item(1493, 15).Value = 0
$ws.Cells.Item(1493, 16).Value = 0
$ws.Cells.Item(1493, 17).Value = 0

# Row 1494
$ws.Cells.Item(1494, 1).NumberFormat = $dateFmt
$ws.Cells.Item(1494, 1).Value = 45509
$ws.Cells.Item(1494, 2).Value = 1472
$ws.Cells.Item(1494, 3).Value = 1589.650024414062
$ws.Cells.Item(1494, 4).Value = 1472
$ws.Cells.Item(1494, 5).Value = 1574.75
$ws.Cells.Item(1494, 6).Value = 1574.75
$ws.Cells.Item(1494, 7).Value = 9352708
$ws.Cells.Item(1494, 8).Value = 2024
$ws.Cells.Item(1494, 9).Value = 8
$ws.Cells.Item(1494, 10).Value = 5
$ws.Cells.Item(1494, 11).Value = 0
$ws.Cells.Item(1494, 12).Value = 0
$ws.Cells.Item(1494, 13).Value = 0
$ws.Cells.Item(1494, 14).Value = 32
$ws.Cells.Item(1494, 15).Value = 0
$ws.Cells.Item(1494, 16).Value = 0
$ws.Cells.Item(1494, 17).Value = 0

# Row 1495
$ws.Cells.Item(1495, 1).NumberFormat = $dateFmt
$ws.Cells.Item(1495, 1).Value = 45516
$ws.Cells.Item(1495, 2).Value = 1574.699951171875
$ws.Cells.Item(1495, 3).Value = 1606.699951171875
$ws.Cells.Item(1495, 4).Value = 1555
$ws.Cells.Item(1495, 5).Value = 1576.099975585938
$ws.Cells.Item(1495, 6).Value = 1576.099975585938
$ws.Cells.Item(1495, 7).Value = 5819526
$ws.Cells.Item(1495, 8).Value = 2024
$ws.Cells.Item(1495, 9).Value = 8
$ws.Cells.Item(1495, 10).Value = 12
$ws.Cells.Item(1495, 11).Value = 0
$ws.Cells.Item(1495, 12).Value = 0
$ws.Cells.Item(1495, 13).Value = 0
$ws.Cells.Item(1495, 14).Value = 33
$ws.Cells.Item(1495, 15).Value = 0
$ws.Cells.Item(1495, 16).Value = 0
$ws.Cells.Item(1495, 17).Value = 0

# Row 1496
$ws.Cells.Item(1496, 1).NumberFormat = $dateFmt
$ws.Cells.Item(1496, 1).Value = 45523
$ws.Cells.Item(1496, 2).Value = 1578.449951171875
$ws.Cells.Item(1496, 3).Value = 1599
$ws.Cells.Item(1496, 4).Value = 1558.199951171875
$ws.Cells.Item(1496, 5).Value = 1574.550048828125
$ws.Cells.Item(1496, 6).Value = 1574.550048828125
$ws.Cells.Item(1496, 7).Value = 4908777
$ws.Cells.Item(1496, 8).Value = 2024
$ws.Cells.Item(1496, 9).Value = 8
$ws.Cells.Item(1496, 10).Value = 19
$ws.Cells.Item(1496, 11).Value = 0
$ws.Cells.Item(1496, 12).Value = 0
$ws.Cells.Item(1496, 13).Value = 0
$ws.Cells.Item(1496, 14).Value = 34
$ws.Cells.Item(1496, 15).Value = 0
$ws.Cells.Item(1496, 16).Value = 0
$ws.Cells.Item(1496, 17).Value = 0

# Row 1497
$ws.Cells.Item(1497, 1).NumberFormat = $dateFmt
$ws.Cells.Item(1497, 1).Value = 45530
$ws.Cells.Item(1497, 2).Value = 1570
$ws.Cells.Item(1497, 3).Value = 1663.699951171875
$ws.Cells.Item(1497, 4).Value = 1566.550048828125
$ws.Cells.Item(1497, 5).Value = 1654.900024414062
$ws.Cells.Item(1497, 6).Value = 1654.900024414062
$ws.Cells.Item(1497, 7).Value = 9708448
$ws.Cells.Item(1497, 8).Value = 2024
$ws.Cells.Item(1497, 9).Value = 8
$ws.Cells.Item(1497, 10).Value = 26
$ws.Cells.Item(1497, 11).Value = 0
$ws.Cells.Item(1497, 12).Value = 0
$ws.Cells.Item(1497, 13).Value = 0
$ws.Cells.Item(1497, 14).Value = 35
$ws.Cells.Item(1497, 15).Value = 0
$ws.Cells.Item(1497, 16).Value = 0
$ws.Cells.Item(1497, 17).Value = 0

# Row 1498
$ws.Cells.Item(1498, 1).NumberFormat = $dateFmt
$ws.Cells.Item(1498, 1).Value = 45537
$ws.Cells.Item(1498, 2).Value = 1669
$ws.Cells.Item(1498, 3).Value = 1681.599975585938
$ws.Cells.Item(1498, 4).Value = 1608.050048828125
$ws.Cells.Item(1498, 5).Value = 1611.050048828125
$ws.Cells.Item(1498, 6).Value = 1611.050048828125
$ws.Cells.Item(1498, 7).Value = 5585820
$ws.Cells.Item(1498, 8).Value = 2024
$ws.Cells.Item(1498, 9).Value = 9
$ws.Cells.Item(1498, 10).Value = 2
$ws.Cells.Item(1498, 11).Value = 0
$ws.Cells.Item(1498, 12).Value = 0
$ws.Cells.Item(1498, 13).Value = 0
$ws.Cells.Item(1498, 14).Value = 36
$ws.Cells.Item(1498, 15).Value = 0
$ws.Cells.Item(1498, 16).Value = 0
$ws.Cells.Item(1498, 17).Value = 0

# Row 1499
$ws.Cells.Item(1499, 1).NumberFormat = $dateFmt
$ws.Cells.Item(1499, 1).Value = 45544
$ws.Cells.Item(1499, 2).Value = 1611.050048828125
$ws.Cells.Item(1499, 3).Value = 1669.199951171875
$ws.Cells.Item(1499, 4).Value = 1610
$ws.Cells.Item(1499, 5).Value = 1659.699951171875
$ws.Cells.Item(1499, 6).Value = 1659.699951171875
$ws.Cells.Item(1499, 7).Value = 5970867
$ws.Cells.Item(1499, 8).Value = 2024
$ws.Cells.Item(1499, 9).Value = 9
$ws.Cells.Item(1499, 10).Value = 9
$ws.Cells.Item(1499, 11).Value = 0
$ws.Cells.Item(1499, 12).Value = 0
$ws.Cells.Item(1499, 13).Value = 0
$ws.Cells.Item(1499, 14).Value = 37
$ws.Cells.Item(1499, 15).Value = 0
$ws.Cells.Item(1499, 16).Value = 0
$ws.Cells.Item(1499, 17).Value = 0

# Row 1500
$ws.Cells.Item(1500, 1).NumberFormat = $dateFmt
$ws.Cells.Item(1500, 1).Value = 45551
$ws.Cells.Item(1500, 2).Value = 1661.199951171875
$ws.Cells.Item(1500, 3).Value = 1683
$ws.Cells.Item(1500, 4).Value = 1609.800048828125
$ws.Cells.Item(1500, 5).Value = 1638.650024414062
$ws.Cells.Item(1500, 6).Value = 1638.650024414062
$ws.Cells.Item(1500, 7).Value = 7278832
$ws.Cells.Item(1500, 8).Value = 2024
$ws.Cells.Item(1500, 9).Value = 9
$ws.Cells.Item(1500, 10).Value = 16
$ws.Cells.Item(1500, 11).Value = 0
$ws.Cells.Item(1500, 12).Value = 0
$ws.Cells.Item(1500, 13).Value = 0
$ws.Cells.Item(1500, 14).Value = 38
$ws.Cells.Item(1500, 15).Value = 0
$ws.Cells.Item(1500, 16).Value = 0
$ws.Cells.Item(1500, 17).Value = 0

# Row 1501
$ws.Cells.Item(1501, 1).NumberFormat = $dateFmt
$ws.Cells.Item(1501, 1).Value = 45558
$ws.Cells.Item(1501, 2).Value = 1646.849975585938
$ws.Cells.Item(1501, 3).Value = 1679.050048828125
$ws.Cells.Item(1501, 4).Value = 1602.300048828125
$ws.Cells.Item(1501, 5).Value = 1672.5
$ws.Cells.Item(1501, 6).Value = 1672.5
$ws.Cells.Item(1501, 7).Value = 9845516
$ws.Cells.Item(1501, 8).Value = 2024
$ws.Cells.Item(1501, 9).Value = 9
$ws.Cells.Item(1501, 10).Value = 23
$ws.Cells.Item(1501, 11).Value = 0
$ws.Cells.Item(1501, 12).Value = 0
$ws.Cells.Item(1501, 13).Value = 0
$ws.Cells.Item(1501, 14).Value = 39
$ws.Cells.Item(1501, 15).Value = 0
$ws.Cells.Item(1501, 16).Value = 0
$ws.Cells.Item(1501, 17).Value = 0

# Row 1502
$ws.Cells.Item(1502, 1).NumberFormat = $dateFmt
$ws.Cells.Item(1502, 1).Value = 45565
$ws.Cells.Item(1502, 2).Value = 1675
$ws.Cells.Item(1502, 3).Value = 1678.900024414062
$ws.Cells.Item(1502, 4).Value = 1577.300048828125
$ws.Cells.Item(1502, 5).Value = 1623.300048828125
$ws.Cells.Item(1502, 6).Value = 1623.300048828125
$ws.Cells.Item(1502, 7).Value = 6163869
$ws.Cells.Item(1502, 8).Value = 2024
$ws.Cells.Item(1502, 9).Value = 9
$ws.Cells.Item(1502, 10).Value = 30
$ws.Cells.Item(1502, 11).Value = 0
$ws.Cells.Item(1502, 12).Value = 0
$ws.Cells.Item(1502, 13).Value = 0
$ws.Cells.Item(1502, 14).Value = 40
$ws.Cells.Item(1502, 15).Value = 0
$ws.Cells.Item(1502, 16).Value = 0
$ws.Cells.Item(1502, 17).Value = 0

# Row 1503
$ws.Cells.Item(1503, 1).NumberFormat = $dateFmt
$ws.Cells.Item(1503, 1).Value = 45572
$ws.Cells.Item(1503, 2).Value = 1623
$ws.Cells.Item(1503, 3).Value = 1702.050048828125
$ws.Cells.Item(1503, 4).Value = 1573.75
$ws.Cells.Item(1503, 5).Value = 1595.75
$ws.Cells.Item(1503, 6).Value = 1595.75
$ws.Cells.Item(1503, 7).Value = 8647323
$ws.Cells.Item(1503, 8).Value = 2024
$ws.Cells.Item(1503, 9).Value = 10
$ws.Cells.Item(1503, 10).Value = 7
$ws.Cells.Item(1503, 11).Value = 0
$ws.Cells.Item(1503, 12).Value = 0
$ws.Cells.Item(1503, 13).Value = 0
$ws.Cells.Item(1503, 14).Value = 41
$ws.Cells.Item(1503, 15).Value = 1
$ws.Cells.Item(1503, 16).Value = 0
$ws.Cells.Item(1503, 17).Value = 0

# Row 1504
$ws.Cells.Item(1504, 1).NumberFormat = $dateFmt
$ws.Cells.Item(1504, 1).Value = 45579
$ws.Cells.Item(1504, 2).Value = 1600
$ws.Cells.Item(1504, 3).Value = 1605.949951171875
$ws.Cells.Item(1504, 4).Value = 1536.349975585938
$ws.Cells.Item(1504, 5).Value = 1551.699951171875
$ws.Cells.Item(1504, 6).Value = 1551.699951171875
$ws.Cells.Item(1504, 7).Value = 10404508
$ws.Cells.Item(1504, 8).Value = 2024
$ws.Cells.Item(1504, 9).Value = 10
$ws.Cells.Item(1504, 10).Value = 14
$ws.Cells.Item(1504, 11).Value = 0
$ws.Cells.Item(1504, 12).Value = 0
$ws.Cells.Item(1504, 13).Value = 0
$ws.Cells.Item(1504, 14).Value = 42
$ws.Cells.Item(1504, 15).Value = 0
$ws.Cells.Item(1504, 16).Value = 0
$ws.Cells.Item(1504, 17).Value = 0

# Row 1505
$ws.Cells.Item(1505, 1).NumberFormat = $dateFmt
$ws.Cells.Item(1505, 1).Value = 45586
$ws.Cells.Item(1505, 2).Value = 1552
$ws.Cells.Item(1505, 3).Value = 1561.25
$ws.Cells.Item(1505, 4).Value = 1478.800048828125
$ws.Cells.Item(1505, 5).Value = 1488.900024414062
$ws.Cells.Item(1505, 6).Value = 1488.900024414062
$ws.Cells.Item(1505, 7).Value = 7334426
$ws.Cells.Item(1505, 8).Value = 2024
$ws.Cells.Item(1505, 9).Value = 10
$ws.Cells.Item(1505, 10).Value = 21
$ws.Cells.Item(1505, 11).Value = 0
$ws.Cells.Item(1505, 12).Value = 0
$ws.Cells.Item(1505, 13).Value = 0
$ws.Cells.Item(1505, 14).Value = 43
$ws.Cells.Item(1505, 15).Value = 0
$ws.Cells.Item(1505, 16).Value = 0
$ws.Cells.Item(1505, 17).Value = 0

# Row 1506
$ws.Cells.Item(1506, 1).NumberFormat = $dateFmt
$ws.Cells.Item(1506, 1).Value = 45593
$ws.Cells.Item(1506, 2).Value = 1496.800048828125
$ws.Cells.Item(1506, 3).Value = 1573.449951171875
$ws.Cells.Item(1506, 4).Value = 1403
$ws.Cells.Item(1506, 5).Value = 1559.550048828125
$ws.Cells.Item(1506, 6).Value = 1559.550048828125
$ws.Cells.Item(1506, 7).Value = 27454663
$ws.Cells.Item(1506, 8).Value = 2024
$ws.Cells.Item(1506, 9).Value = 10
$ws.Cells.Item(1506, 10).Value = 28
$ws.Cells.Item(1506, 11).Value = 0
$ws.Cells.Item(1506, 12).Value = 0
$ws.Cells.Item(1506, 13).Value = 0
$ws.Cells.Item(1506, 14).Value = 44
$ws.Cells.Item(1506, 15).Value = 0
$ws.Cells.Item(1506, 16).Value = 0
$ws.Cells.Item(1506, 17).Value = 0

# Row 1507
$ws.Cells.Item(1507, 1).NumberFormat = $dateFmt
$ws.Cells.Item(1507, 1).Value = 45600
$ws.Cells.Item(1507, 2).Value = 1584
$ws.Cells.Item(1507, 3).Value = 1612.349975585938
$ws.Cells.Item(1507, 4).Value = 1558.550048828125
$ws.Cells.Item(1507, 5).Value = 1592.599975585938
$ws.Cells.Item(1507, 6).Value = 1592.599975585938
$ws.Cells.Item(1507, 7).Value = 11734078
$ws.Cells.Item(1507, 8).Value = 2024
$ws.Cells.Item(1507, 9).Value = 11
$ws.Cells.Item(1507, 10).Value = 4
$ws.Cells.Item(1507, 11).Value = 0
$ws.Cells.Item(1507, 12).Value = 0
$ws.Cells.Item(1507, 13).Value = 0
$ws.Cells.Item(1507, 14).Value = 45
$ws.Cells.Item(1507, 15).Value = 0
$ws.Cells.Item(1507, 16).Value = 0
$ws.Cells.Item(1507, 17).Value = 0

# Row 1508
$ws.Cells.Item(1508, 1).NumberFormat = $dateFmt
$ws.Cells.Item(1508, 1).Value = 45607
$ws.Cells.Item(1508, 2).Value = 1590.099975585938
$ws.Cells.Item(1508, 3).Value = 1599.75
$ws.Cells.Item(1508, 4).Value = 1490.25
$ws.Cells.Item(1508, 5).Value = 1499.75
$ws.Cells.Item(1508, 6).Value = 1499.75
$ws.Cells.Item(1508, 7).Value = 8998865
$ws.Cells.Item(1508, 8).Value = 2024
$ws.Cells.Item(1508, 9).Value = 11
$ws.Cells.Item(1508, 10).Value = 11
$ws.Cells.Item(1508, 11).Value = 0
$ws.Cells.Item(1508, 12).Value = 0
$ws.Cells.Item(1508, 13).Value = 0
$ws.Cells.Item(1508, 14).Value = 46
$ws.Cells.Item(1508, 15).Value = 0
$ws.Cells.Item(1508, 16).Value = 0
$ws.Cells.Item(1508, 17).Value = 0
